$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "юрист"
$ws.Range("D3").Value = "садовник"
$ws.Range("D3").Select()
